$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Fecha actual" marker (D13 -> D14): move the highlighted "fecha actual"
# flag down to the row for Thursday 17 August (class date), removing it
# from the previous class date.
$ws.Range("D13").Copy($ws.Range("D14"))
$ws.Range("D13").Clear()

# Mark the 17 Aug class as reviewed/ok.
$ws.Range("G14").Value = "ok"

# New notes at the bottom of the sheet.
$ws.Range("B23").Value = "sugerencias: entregar prueba antes del 2/9 "
$ws.Range("B24").Value = "fecha entrega maxima: 4/5? "

# Leave the selection where the author left off.
$ws.Range("B25").Select() | Out-Null
